# Przygotowanie (wszystkich) operacji dla modelu Country(Panstwo)
# Insert a new day of data at row 2, shifting existing rows down, then
# add the new day's data at row 4 (date 43440 / 2018-12-06), and refresh
# the helper formulas for rows 3-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dates: shift the tracked dates back by one day (a new earliest
#     date, 43438, becomes the top row; what used to be row 3's date
#     moves to row 3 still but with the new, one-day-earlier value) ---
$ws.Range("A2").Value = 43438
$ws.Range("A3").Value = 43439

# --- New day of data added at row 4 ---
$ws.Range("A4").Value = 43440
$ws.Range("B4").Value = 88603

# --- Row 3: C3 stays B3/B2 (unchanged); D3 now compares B3 directly
#     against the baseline B2 instead of going through C3 ---
$ws.Range("C3").Formula = "=B3/B2"
$ws.Range("D3").Formula = "=B3/`$B`$2"

# --- Row 4: day-over-day ratio (shared group anchored here) and the
#     direct-to-baseline ratio ---
$ws.Range("C4").Formula = "=B4/B3"
$ws.Range("D4").Formula = "=B4/`$B`$2"

# --- Rows 5-8: shared formula groups re-anchored at row 5 now that
#     row 4 holds real data ---
$ws.Range("C5").Formula = "=B5/B4"
$ws.Range("D5").Formula = "=B5/`$B`$2"
$ws.Range("C6").Formula = "=B6/B5"
$ws.Range("D6").Formula = "=B6/`$B`$2"
$ws.Range("C7").Formula = "=B7/B6"
$ws.Range("D7").Formula = "=B7/`$B`$2"
$ws.Range("C8").Formula = "=B8/B7"
$ws.Range("D8").Formula = "=B8/`$B`$2"

# --- Move the active selection to G5 ---
$ws.Range("G5").Select()
